# Fruta / hortaliza, semanal
# Insert 3 new rows (weekly update) before the former row 49, shifting the
# existing rows 49-113 down to 52-116, then populate the new rows with the
# new "Royal Glory" price entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 49 (pushes old rows 49-113 down to 52-116)
$ws.Range("A49:T51").EntireRow.Insert()

# New row 49: Royal Glory / Especial
$ws.Range("A49").Value = 11
$ws.Range("B49").Value = "Vega Monumental Concepción"
$ws.Range("C49").Value = "Bíobío"
$ws.Range("D49").Value = 44546
$ws.Range("E49").Value = 8
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100103
$ws.Range("H49").Value = "Frutos de hueso (carozo)"
$ws.Range("I49").Value = 100103004
$ws.Range("J49").Value = "Durazno"
$ws.Range("K49").Value = "Royal Glory"
$ws.Range("L49").Value = "Especial"
$ws.Range("M49").Value = 50
$ws.Range("N49").Value = 17000
$ws.Range("O49").Value = 17000
$ws.Range("P49").Value = 17000
$ws.Range("Q49").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R49").Value = "Región de O'Higgins"
$ws.Range("S49").Value = 1062
$ws.Range("T49").Value = 16

# New row 50: Royal Glory / Primera
$ws.Range("A50").Value = 11
$ws.Range("B50").Value = "Vega Monumental Concepción"
$ws.Range("C50").Value = "Bíobío"
$ws.Range("D50").Value = 44546
$ws.Range("E50").Value = 8
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100103
$ws.Range("H50").Value = "Frutos de hueso (carozo)"
$ws.Range("I50").Value = 100103004
$ws.Range("J50").Value = "Durazno"
$ws.Range("K50").Value = "Royal Glory"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 50
$ws.Range("N50").Value = 15000
$ws.Range("O50").Value = 15000
$ws.Range("P50").Value = 15000
$ws.Range("Q50").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R50").Value = "Región de O'Higgins"
$ws.Range("S50").Value = 938
$ws.Range("T50").Value = 16

# New row 51: Royal Glory / Segunda
$ws.Range("A51").Value = 11
$ws.Range("B51").Value = "Vega Monumental Concepción"
$ws.Range("C51").Value = "Bíobío"
$ws.Range("D51").Value = 44546
$ws.Range("E51").Value = 8
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100103
$ws.Range("H51").Value = "Frutos de hueso (carozo)"
$ws.Range("I51").Value = 100103004
$ws.Range("J51").Value = "Durazno"
$ws.Range("K51").Value = "Royal Glory"
$ws.Range("L51").Value = "Segunda"
$ws.Range("M51").Value = 50
$ws.Range("N51").Value = 13000
$ws.Range("O51").Value = 13000
$ws.Range("P51").Value = 13000
$ws.Range("Q51").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R51").Value = "Región de O'Higgins"
$ws.Range("S51").Value = 812
$ws.Range("T51").Value = 16
